$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, 5, 3),
    @(2, 7, 11.12602466666667),
    @(2, 8, 33.378074),
    @(2, 9, 0.7147750296730129),
    @(2, 10, 0.714775029673013),
    @(2, 11, 3),
    @(2, 13, 42.23928166666666),
    @(2, 14, 126.717845),
    @(2, 15, 0.7297675404946526),
    @(2, 16, 0.7297675404946528),
    @(2, 17, 469.9552897256144),
    @(2, 18, 4229.59760753053),
    @(2, 19, 0.521619615411467),
    @(2, 20, 0.5216196154114672),
    @(3, 5, 3),
    @(3, 7, 11.12602466666667),
    @(3, 8, 33.378074),
    @(3, 9, 0.7147750296730129),
    @(3, 10, 0.714775029673013),
    @(3, 11, 3),
    @(3, 13, 4.624984),
    @(3, 14, 13.874952),
    @(3, 15, 0.07990579066051323),
    @(3, 16, 0.07990579066051323),
    @(3, 17, 51.45768606693867),
    @(3, 18, 463.119174602448),
    @(3, 19, 0.0571146638904139),
    @(3, 20, 0.05711466389041391),
    @(4, 5, 3),
    @(4, 7, 11.12602466666667),
    @(4, 8, 33.378074),
    @(4, 9, 0.7147750296730129),
    @(4, 10, 0.714775029673013),
    @(4, 11, 3),
    @(4, 13, 2.151907),
    @(4, 14, 6.455721),
    @(4, 15, 0.03717847029587412),
    @(4, 16, 0.03717847029587412),
    @(4, 17, 23.94217036237266),
    @(4, 18, 215.479533261354),
    @(4, 19, 0.02657424220893065),
    @(4, 20, 0.02657424220893065),
    @(5, 5, 3),
    @(5, 7, 11.12602466666667),
    @(5, 8, 33.378074),
    @(5, 9, 0.7147750296730129),
    @(5, 10, 0.714775029673013),
    @(5, 11, 3),
    @(5, 13, 8.864288333333333),
    @(5, 14, 26.592865),
    @(5, 15, 0.1531481985489599),
    @(5, 16, 0.15314819854896),
    @(5, 17, 98.62429064911221),
    @(5, 18, 887.6186158420099),
    @(5, 19, 0.1094665081622013),
    @(5, 20, 0.1094665081622014),
    @(6, 5, 3),
    @(6, 7, 2.698866333333334),
    @(6, 8, 8.096599000000001),
    @(6, 9, 0.1733846833246127),
    @(6, 10, 0.1733846833246127),
    @(6, 11, 3),
    @(6, 13, 42.23928166666666),
    @(6, 14, 126.717845),
    @(6, 15, 0.7297675404946526),
    @(6, 16, 0.7297675404946528),
    @(6, 17, 113.9981752343506),
    @(6, 18, 1025.983577109155),
    @(6, 19, 0.1265305139092468),
    @(6, 20, 0.1265305139092468),
    @(7, 5, 3),
    @(7, 7, 2.698866333333334),
    @(7, 8, 8.096599000000001),
    @(7, 9, 0.1733846833246127),
    @(7, 10, 0.1733846833246127),
    @(7, 11, 3),
    @(7, 13, 4.624984),
    @(7, 14, 13.874952),
    @(7, 15, 0.07990579066051323),
    @(7, 16, 0.07990579066051323),
    @(7, 17, 12.48221360980533),
    @(7, 18, 112.339922488248),
    @(7, 19, 0.01385444020947588),
    @(7, 20, 0.01385444020947588),
    @(8, 5, 3),
    @(8, 7, 2.698866333333334),
    @(8, 8, 8.096599000000001),
    @(8, 9, 0.1733846833246127),
    @(8, 10, 0.1733846833246127),
    @(8, 11, 3),
    @(8, 13, 2.151907),
    @(8, 14, 6.455721),
    @(8, 15, 0.03717847029587412),
    @(8, 16, 0.03717847029587412),
    @(8, 17, 5.807709354764334),
    @(8, 18, 52.26938419287901),
    @(8, 19, 0.006446177298743651),
    @(8, 20, 0.006446177298743652),
    @(9, 5, 3),
    @(9, 7, 2.698866333333334),
    @(9, 8, 8.096599000000001),
    @(9, 9, 0.1733846833246127),
    @(9, 10, 0.1733846833246127),
    @(9, 11, 3),
    @(9, 13, 8.864288333333333),
    @(9, 14, 26.592865),
    @(9, 15, 0.1531481985489599),
    @(9, 16, 0.15314819854896),
    @(9, 17, 23.92352935179278),
    @(9, 18, 215.311764166135),
    @(9, 19, 0.02655355190714632),
    @(9, 20, 0.02655355190714633),
    @(10, 5, 3),
    @(10, 7, 1.076470666666667),
    @(10, 8, 3.229412),
    @(10, 9, 0.06915626881666041),
    @(10, 10, 0.06915626881666043),
    @(10, 11, 3),
    @(10, 13, 42.23928166666666),
    @(10, 14, 126.717845),
    @(10, 15, 0.7297675404946526),
    @(10, 16, 0.7297675404946528),
    @(10, 17, 45.46934769523777),
    @(10, 18, 409.22412925714),
    @(10, 19, 0.05046800020412131),
    @(10, 20, 0.05046800020412133),
    @(11, 5, 3),
    @(11, 7, 1.076470666666667),
    @(11, 8, 3.229412),
    @(11, 9, 0.06915626881666041),
    @(11, 10, 0.06915626881666043),
    @(11, 11, 3),
    @(11, 13, 4.624984),
    @(11, 14, 13.874952),
    @(11, 15, 0.07990579066051323),
    @(11, 16, 0.07990579066051323),
    @(11, 17, 4.978659609802667),
    @(11, 18, 44.807936488224),
    @(11, 19, 0.005525986338926246),
    @(11, 20, 0.005525986338926247),
    @(12, 5, 3),
    @(12, 7, 1.076470666666667),
    @(12, 8, 3.229412),
    @(12, 9, 0.06915626881666041),
    @(12, 10, 0.06915626881666043),
    @(12, 11, 3),
    @(12, 13, 2.151907),
    @(12, 14, 6.455721),
    @(12, 15, 0.03717847029587412),
    @(12, 16, 0.03717847029587412),
    @(12, 17, 2.316464762894666),
    @(12, 18, 20.848182866052),
    @(12, 19, 0.002571124285973695),
    @(12, 20, 0.002571124285973695),
    @(13, 5, 3),
    @(13, 7, 1.076470666666667),
    @(13, 8, 3.229412),
    @(13, 9, 0.06915626881666041),
    @(13, 10, 0.06915626881666043),
    @(13, 11, 3),
    @(13, 13, 8.864288333333333),
    @(13, 14, 26.592865),
    @(13, 15, 0.1531481985489599),
    @(13, 16, 0.15314819854896),
    @(13, 17, 9.542146371708887),
    @(13, 18, 85.87931734538),
    @(13, 19, 0.01059115798763916),
    @(13, 20, 0.01059115798763916),
    @(14, 5, 3),
    @(14, 7, 0.6644096666666667),
    @(14, 8, 1.993229),
    @(14, 9, 0.04268401818571407),
    @(14, 10, 0.04268401818571407),
    @(14, 11, 3),
    @(14, 13, 42.23928166666666),
    @(14, 14, 126.717845),
    @(14, 15, 0.7297675404946526),
    @(14, 16, 0.7297675404946528),
    @(14, 17, 28.06418705238944),
    @(14, 18, 252.577683471505),
    @(14, 19, 0.03114941096981758),
    @(14, 20, 0.03114941096981759),
    @(15, 5, 3),
    @(15, 7, 0.6644096666666667),
    @(15, 8, 1.993229),
    @(15, 9, 0.04268401818571407),
    @(15, 10, 0.04268401818571407),
    @(15, 11, 3),
    @(15, 13, 4.624984),
    @(15, 14, 13.874952),
    @(15, 15, 0.07990579066051323),
    @(15, 16, 0.07990579066051323),
    @(15, 17, 3.072884077778667),
    @(15, 18, 27.655956700008),
    @(15, 19, 0.003410700221697208),
    @(15, 20, 0.003410700221697208),
    @(16, 5, 3),
    @(16, 7, 0.6644096666666667),
    @(16, 8, 1.993229),
    @(16, 9, 0.04268401818571407),
    @(16, 10, 0.04268401818571407),
    @(16, 11, 3),
    @(16, 13, 2.151907),
    @(16, 14, 6.455721),
    @(16, 15, 0.03717847029587412),
    @(16, 16, 0.03717847029587412),
    @(16, 17, 1.429747812567667),
    @(16, 18, 12.867730313109),
    @(16, 19, 0.001586926502226121),
    @(16, 20, 0.001586926502226121),
    @(17, 5, 3),
    @(17, 7, 0.6644096666666667),
    @(17, 8, 1.993229),
    @(17, 9, 0.04268401818571407),
    @(17, 10, 0.04268401818571407),
    @(17, 11, 3),
    @(17, 13, 8.864288333333333),
    @(17, 14, 26.592865),
    @(17, 15, 0.1531481985489599),
    @(17, 16, 0.15314819854896),
    @(17, 17, 5.889518856787221),
    @(17, 18, 53.00566971108501),
    @(17, 19, 0.006536980491973155),
    @(17, 20, 0.006536980491973157),
)

foreach ($u in $updates) {
    $r = $u[0]
    $c = $u[1]
    $val = $u[2]
    $ws.Cells.Item($r, $c).Value2 = $val
}

Write-Host "Done applying updates: " $updates.Count